$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.718.04"
$ws.Range("E2").Value = "  -4.02%  "

$ws.Range("D3").Value = "1.814.77"
$ws.Range("E3").Value = "  -3.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "278.25"
$ws.Range("E5").Value = "  -7.57%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5084"
$ws.Range("E7").Value = "  -4.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3531"
$ws.Range("E8").Value = "  -5.67%  "

$ws.Range("E9").Value = "  -2.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06678"
$ws.Range("E10").Value = "  -7.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.89"
$ws.Range("E11").Value = "  -8.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8216"
$ws.Range("E12").Value = "  -7.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07858"
$ws.Range("E13").Value = "  -3.76%  "

$ws.Range("D14").Value = "1.824.30"
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.070"
$ws.Range("E15").Value = "  -4.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.60"
$ws.Range("E16").Value = "  -5.71%  "

$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.08"
$ws.Range("E18").Value = "  -5.14%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008030"
$ws.Range("E19").Value = "  -5.77%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").Value = "25.767.14"
$ws.Range("E21").Value = "  -3.97%  "

$ws.Range("E22").Value = "  -4.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.983"
$ws.Range("E23").Value = "  -6.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.107"
$ws.Range("E24").Value = "  -4.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.240"
$ws.Range("E25").Value = "  -3.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.52"
$ws.Range("E26").Value = "  -2.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.666"
$ws.Range("E27").Value = "  -3.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.12"
$ws.Range("E28").Value = "  -5.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "109.00"
$ws.Range("E29").Value = "  -4.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.329"
$ws.Range("E30").Value = "  -8.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.218"
$ws.Range("E31").Value = "  -8.90%  "

$ws.Range("E32").Value = "  -4.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04863"
$ws.Range("E33").Value = "  -3.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7269"
$ws.Range("E34").Value = "  -9.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.886"
$ws.Range("E35").Value = "  -1.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.130"
$ws.Range("E36").Value = "  -3.93%  "

$ws.Range("E37").Value = "  -1.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.365"
$ws.Range("E38").Value = "  -12.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01849"
$ws.Range("E39").Value = "  -5.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5133"
$ws.Range("E40").Value = "  -16.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9688"
$ws.Range("E41").Value = "  -8.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "114.22"
$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.223"
$ws.Range("E43").Value = "  -4.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.007"
$ws.Range("E44").Value = "  -8.73%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4524"
$ws.Range("E46").Value = "  -13.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1367"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.37"
$ws.Range("E48").Value = "  -3.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.154"
$ws.Range("E49").Value = "  -8.17%  "

$ws.Range("E50").Value = "  -8.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05821"
$ws.Range("E51").Value = "  -3.84%  "
